# Generate Report for Handback
# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps for the 7b9c8ed8-... row on the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E3").Value = "2016-03-13 08:45:26"
$zhcn.Range("H3").Value = "2016-03-13 08:45:44"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E3").Value = "2016-03-13 08:45:32"
$dede.Range("H3").Value = "2016-03-13 08:45:50"
